$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.157.60"
$ws.Range("E2").Value = "  +5.78%  "
$ws.Range("D3").Value = "3.744.97"
$ws.Range("E3").Value = "  +20.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.27"
$ws.Range("E5").Value = "  +7.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.94"
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("D7").Value = "3.735.75"
$ws.Range("E7").Value = "  +20.31%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +5.76%  "
$ws.Range("E10").Value = "  +7.84%  "
$ws.Range("E11").Value = "  +3.22%  "
$ws.Range("E12").Value = "  +7.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.56"
$ws.Range("E13").Value = "  +11.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000257"
$ws.Range("E14").Value = "  +6.14%  "
$ws.Range("D15").Value = "4.362.43"
$ws.Range("E15").Value = "  +20.33%  "
$ws.Range("D16").Value = "3.732.68"
$ws.Range("E16").Value = "  +20.19%  "
$ws.Range("D17").Value = "71.214.64"
$ws.Range("E17").Value = "  +6.02%  "
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.55"
$ws.Range("E19").Value = "  +7.07%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.95"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "519.77"
$ws.Range("E21").Value = "  +5.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.39"
$ws.Range("E22").Value = "  +21.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.747"
$ws.Range("E23").Value = "  +8.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.54"
$ws.Range("E24").Value = "  +11.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.70"
$ws.Range("E25").Value = "  +6.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.58"
$ws.Range("E26").Value = "  +7.67%  "
$ws.Range("E27").Value = "  +10.29%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +10.00%  "
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.90"
$ws.Range("E31").Value = "  +11.56%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000113"
$ws.Range("E32").Value = "  +19.50%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.20"
$ws.Range("E33").Value = "  +14.34%  "
$ws.Range("E34").Value = "  +4.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  +10.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.15"
$ws.Range("E37").Value = "  +10.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.346"
$ws.Range("E38").Value = "  +10.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.24"
$ws.Range("E39").Value = "  +10.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.135"
$ws.Range("E40").Value = "  +9.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.74"
$ws.Range("E41").Value = "  +5.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "435.81"
$ws.Range("E42").Value = "  +17.06%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.163.41"
$ws.Range("E43").Value = "  +12.86%  "
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.73"
$ws.Range("E44").Value = "  -5.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.86"
$ws.Range("E45").Value = "  +6.72%  "
$ws.Range("E46").Value = "  +5.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0366"
$ws.Range("E47").Value = "  +5.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.16"
$ws.Range("E48").Value = "  +9.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.93"
$ws.Range("E49").Value = "  +3.80%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.49"
$ws.Range("E51").Value = "  +8.62%  "
